$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- original row 21
$ws.Range("D2").Value = 44355
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 18000
$ws.Range("O2").Value = 18000
$ws.Range("P2").Value = 18000
$ws.Range("Q2").Value = "$/caja 18 kilos granel"
$ws.Range("R2").Value = "Provincia de Limarí"
$ws.Range("S2").Value = 1000
$ws.Range("T2").Value = 18

# Row 3 <- original row 8
$ws.Range("D3").Value = 44299
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("Q3").Value = "$/caja 15 kilos granel"
$ws.Range("R3").Value = "Provincia de Curicó"
$ws.Range("S3").Value = 1000
$ws.Range("T3").Value = 15

# Row 4: unchanged (maps to itself)
# Row 5: unchanged (maps to itself)
# Row 6 <- original row 9
$ws.Range("D6").Value = 44291
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 150
$ws.Range("N6").Value = 12000
$ws.Range("O6").Value = 12000
$ws.Range("P6").Value = 12000
$ws.Range("Q6").Value = "$/caja 15 kilos granel"
$ws.Range("R6").Value = "Región Metropolitana"
$ws.Range("S6").Value = 800
$ws.Range("T6").Value = 15

# Row 7 <- original row 23
$ws.Range("D7").Value = 45099
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 17000
$ws.Range("O7").Value = 17000
$ws.Range("P7").Value = 17000
$ws.Range("Q7").Value = "$/caja 18 kilos granel"
$ws.Range("R7").Value = "Provincia de Limarí"
$ws.Range("S7").Value = 944
$ws.Range("T7").Value = 18

# Row 8 <- original row 13
$ws.Range("D8").Value = 44326
$ws.Range("L8").Value = "Especial"
$ws.Range("M8").Value = 300
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 20000
$ws.Range("P8").Value = 20000
$ws.Range("Q8").Value = "$/caja 18 kilos granel"
$ws.Range("R8").Value = "Provincia de Limarí"
$ws.Range("S8").Value = 1111
$ws.Range("T8").Value = 18

# Row 9 <- original row 22
$ws.Range("D9").Value = 45062
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 200
$ws.Range("N9").Value = 15000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 15000
$ws.Range("Q9").Value = "$/caja 15 kilos granel"
$ws.Range("R9").Value = "Provincia de Curicó"
$ws.Range("S9").Value = 1000
$ws.Range("T9").Value = 15

# Row 10 <- original row 15
$ws.Range("D10").Value = 44692
$ws.Range("L10").Value = "Especial"
$ws.Range("M10").Value = 150
$ws.Range("N10").Value = 17000
$ws.Range("O10").Value = 17000
$ws.Range("P10").Value = 17000
$ws.Range("Q10").Value = "$/caja 18 kilos granel"
$ws.Range("R10").Value = "Provincia de Limarí"
$ws.Range("S10").Value = 944
$ws.Range("T10").Value = 18

# Row 11 <- original row 12
$ws.Range("D11").Value = 44340
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 230
$ws.Range("N11").Value = 20000
$ws.Range("O11").Value = 20000
$ws.Range("P11").Value = 20000
$ws.Range("Q11").Value = "$/caja 18 kilos granel"
$ws.Range("R11").Value = "Provincia de Limarí"
$ws.Range("S11").Value = 1111
$ws.Range("T11").Value = 18

# Row 12 <- original row 17
$ws.Range("D12").Value = 44714
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 100
$ws.Range("N12").Value = 20000
$ws.Range("O12").Value = 20000
$ws.Range("P12").Value = 20000
$ws.Range("Q12").Value = "$/caja 18 kilos granel"
$ws.Range("R12").Value = "Provincia de Limarí"
$ws.Range("S12").Value = 1111
$ws.Range("T12").Value = 18

# Row 13 <- original row 11
$ws.Range("D13").Value = 44319
$ws.Range("L13").Value = "Especial"
$ws.Range("M13").Value = 120
$ws.Range("N13").Value = 20000
$ws.Range("O13").Value = 20000
$ws.Range("P13").Value = 20000
$ws.Range("Q13").Value = "$/caja 18 kilos granel"
$ws.Range("R13").Value = "Provincia de Limarí"
$ws.Range("S13").Value = 1111
$ws.Range("T13").Value = 18

# Row 14 <- original row 20
$ws.Range("D14").Value = 44680
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 200
$ws.Range("N14").Value = 15000
$ws.Range("O14").Value = 15000
$ws.Range("P14").Value = 15000
$ws.Range("Q14").Value = "$/caja 15 kilos granel"
$ws.Range("R14").Value = "Provincia de Limarí"
$ws.Range("S14").Value = 1000
$ws.Range("T14").Value = 15

# Row 15 <- original row 18
$ws.Range("D15").Value = 45083
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 120
$ws.Range("N15").Value = 17000
$ws.Range("O15").Value = 17000
$ws.Range("P15").Value = 17000
$ws.Range("Q15").Value = "$/caja 18 kilos granel"
$ws.Range("R15").Value = "Región de O'Higgins"
$ws.Range("S15").Value = 944
$ws.Range("T15").Value = 18

# Row 16 <- original row 14
$ws.Range("D16").Value = 45054
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 200
$ws.Range("N16").Value = 14000
$ws.Range("O16").Value = 14000
$ws.Range("P16").Value = 14000
$ws.Range("Q16").Value = "$/caja 18 kilos granel"
$ws.Range("R16").Value = "Provincia de Curicó"
$ws.Range("S16").Value = 778
$ws.Range("T16").Value = 18

# Row 17 <- original row 19
$ws.Range("D17").Value = 44348
$ws.Range("L17").Value = "Especial"
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = 20000
$ws.Range("O17").Value = 20000
$ws.Range("P17").Value = 20000
$ws.Range("Q17").Value = "$/caja 18 kilos granel"
$ws.Range("R17").Value = "Provincia de Limarí"
$ws.Range("S17").Value = 1111
$ws.Range("T17").Value = 18

# Row 18 <- original row 24
$ws.Range("D18").Value = 45085
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 280
$ws.Range("N18").Value = 17000
$ws.Range("O18").Value = 18000
$ws.Range("P18").Value = 17357
$ws.Range("Q18").Value = "$/caja 18 kilos granel"
$ws.Range("R18").Value = "Provincia de Curicó"
$ws.Range("S18").Value = 964
$ws.Range("T18").Value = 18

# Row 19 <- original row 25
$ws.Range("D19").Value = 44354
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 100
$ws.Range("N19").Value = 18000
$ws.Range("O19").Value = 18000
$ws.Range("P19").Value = 18000
$ws.Range("Q19").Value = "$/caja 18 kilos granel"
$ws.Range("R19").Value = "Provincia de Limarí"
$ws.Range("S19").Value = 1000
$ws.Range("T19").Value = 18

# Row 20 <- original row 7
$ws.Range("D20").Value = 44691
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = 17000
$ws.Range("O20").Value = 17000
$ws.Range("P20").Value = 17000
$ws.Range("Q20").Value = "$/caja 18 kilos granel"
$ws.Range("R20").Value = "Provincia de Limarí"
$ws.Range("S20").Value = 944
$ws.Range("T20").Value = 18

# Row 21 <- original row 16
$ws.Range("D21").Value = 44328
$ws.Range("L21").Value = "Especial"
$ws.Range("M21").Value = 250
$ws.Range("N21").Value = 20000
$ws.Range("O21").Value = 20000
$ws.Range("P21").Value = 20000
$ws.Range("Q21").Value = "$/caja 18 kilos granel"
$ws.Range("R21").Value = "Provincia de Limarí"
$ws.Range("S21").Value = 1111
$ws.Range("T21").Value = 18

# Row 22 <- original row 10
$ws.Range("D22").Value = 44316
$ws.Range("L22").Value = "Especial"
$ws.Range("M22").Value = 300
$ws.Range("N22").Value = 20000
$ws.Range("O22").Value = 20000
$ws.Range("P22").Value = 20000
$ws.Range("Q22").Value = "$/caja 18 kilos granel"
$ws.Range("R22").Value = "Provincia de Limarí"
$ws.Range("S22").Value = 1111
$ws.Range("T22").Value = 18

# Row 23 <- original row 2
$ws.Range("D23").Value = 44294
$ws.Range("L23").Value = "Primera"
$ws.Range("M23").Value = 50
$ws.Range("N23").Value = 12000
$ws.Range("O23").Value = 12000
$ws.Range("P23").Value = 12000
$ws.Range("Q23").Value = "$/caja 15 kilos granel"
$ws.Range("R23").Value = "Región Metropolitana"
$ws.Range("S23").Value = 800
$ws.Range("T23").Value = 15

# Row 24 <- original row 3
$ws.Range("D24").Value = 44342
$ws.Range("L24").Value = "Especial"
$ws.Range("M24").Value = 300
$ws.Range("N24").Value = 20000
$ws.Range("O24").Value = 20000
$ws.Range("P24").Value = 20000
$ws.Range("Q24").Value = "$/caja 18 kilos granel"
$ws.Range("R24").Value = "Provincia de Limarí"
$ws.Range("S24").Value = 1111
$ws.Range("T24").Value = 18

# Row 25 <- original row 6
$ws.Range("D25").Value = 45097
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 200
$ws.Range("N25").Value = 17000
$ws.Range("O25").Value = 17000
$ws.Range("P25").Value = 17000
$ws.Range("Q25").Value = "$/caja 18 kilos granel"
$ws.Range("R25").Value = "Provincia de Limarí"
$ws.Range("S25").Value = 944
$ws.Range("T25").Value = 18
